$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 29 :: 2/20/2020 (serial 43881) ----
$ws.Range("A26").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = 43881

$ws.Range("B28").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("B29").Value = "17:00-19:00 in class"

$ws.Range("C28").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = "N.A."

$ws.Range("D19").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = "Follow the lecture with professor"

$ws.Range("F28").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = "leant to dive deeper within a project from architecture perspective, also social context perspective"

$ws.Range("F28").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = "Instead of doing more research on class level and feature level of a project, we learnt to do more work based on an architectural perspective. This will give us a more concreted sense of the whole program. Also we learnt to know the importance of looking at the social context of a project before deciding to contribute. It might be frustrated if maintainers do not maintain this project any more when you contribute."

$ws.Range("G28").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = "Average"

# ---- Row 30 :: 2/25/2020 (serial 43886) ----
$ws.Range("A26").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A30").Value = 43886

$ws.Range("B28").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$ws.Range("B30").Value = "13:30 - 17:00"

$ws.Range("C28").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = "Soobin"

$ws.Range("D28").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = "finish homework3"

$ws.Range("F28").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = "talked about interesting open issues and PRs, finished social context of our program"

$ws.Range("F28").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = "Some of these PRs are quite ridiculous. Some people they just want to be a contributor by doing minimal, such as fixing typo in commands. And that is all they do. So funny and things also happen in reality."

$ws.Range("G28").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = "Nervous"

# ---- Row 31 :: 2/26/2020 (serial 43887) ----
$ws.Range("A26").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A31").Value = 43887

$ws.Range("B28").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("B31").Value = "14:00 - 17:00 & 19:30 - 21:00"

$ws.Range("C28").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("C31").Value = "Soobin"

$ws.Range("D28").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = "finish homework3"

$ws.Range("F28").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = "finished architecture recovering"

$ws.Range("F28").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F31").Value = "What we did in class by dragging and grouping Pacman classes in UML diagram is kinda a bottom up comprehension. But when it comes to a 100K LOC program, it is impossible to use that strategy since there are SOOOO many classes and interfaces. So we used a top down way to comprehend all the meaningful features first then top down confirm our hypothesis of the architecture by taking a look at these folders’ name, files implementations. It works fine"

$ws.Range("G28").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").Value = "Good, relief"
